$wb = $excel.ActiveWorkbook
$wsParts = $wb.Worksheets.Item("Shield Parts")
$wsAdd   = $wb.Worksheets.Item("To Add")

# ------------------------------------------------------------------
# "Shield Parts" sheet: remove the "Quad Buffer" line item (row 17).
# Leave the shared F/H/J formulas alone; they'll recompute to 0 once
# their inputs are cleared.
# ------------------------------------------------------------------
$wsParts.Range("C17").ClearContents()
$wsParts.Range("D17").ClearContents()
$wsParts.Range("E17").ClearContents()
$wsParts.Range("G17").ClearContents()
$wsParts.Range("I17").ClearContents()
$wsParts.Range("K17").ClearContents()
$wsParts.Range("L17").ClearContents()
$wsParts.Range("M17").ClearContents()
$wsParts.Range("N17").ClearContents()

# The N17 hyperlink (to the now-removed part's datasheet) needs to move
# down to N18, which is now the first row with a link. Rebuild the
# hyperlink collection with the same addresses, shifting N17 -> N18.
[void]$wsParts.Hyperlinks.Delete()

[void]$wsParts.Hyperlinks.Add($wsParts.Range("N9"),  "http://www.mouser.com/ProductDetail/Maxim-Integrated/DS3231SNTR/?qs=sGAEpiMZZMuuBt6TL7D%2f6PgM9QV8pLmA")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N10"), "http://www.mouser.com/ProductDetail/Keystone-Electronics/3000/?qs=sGAEpiMZZMtT9MhkajLHrnU1d13jcSgSROM9zhZkF8A%3d")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N11"), "http://www.mouser.com/Search/ProductDetail.aspx?qs=sGAEpiMZZMtz8P%2feuiupSd2F%2fX%2ffEmeEyLDGD5JMOeY%3d")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N12"), "https://www.sparkfun.com/products/8077")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N15"), "https://www.sparkfun.com/products/102")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N16"), "http://www.newark.com/te-connectivity/2041021-3/memory-card-connector-sd-9-position/dp/35R2925?CMP=AFC-QO1721829242?gross_price=")
[void]$wsParts.Hyperlinks.Add($wsParts.Range("N18"), "http://www.mouser.com/ProductDetail/Texas-Instruments/SN74AHC125MDREP/?qs=sGAEpiMZZMuiiWkaIwCK2S7iisUJKLbkCxHZbiEL4Hk%3d")

# ------------------------------------------------------------------
# "To Add" sheet: new note about a transistor for soil sensors.
# ------------------------------------------------------------------
$wsAdd.Range("A10").Value = "Transistor for soil sensors"

# ------------------------------------------------------------------
# Selection / view bookkeeping to match the saved workbook state.
# ------------------------------------------------------------------
$wsParts.Activate()
$wsParts.Range("B24").Select()

$wsAdd.Activate()
$wsAdd.Range("A11").Select()
